$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F7: GDP Nowcast present value update ---
$ws.Range("F7").Value = 0.3276016540317022

# --- Remove "latest update" highlight (style 48 -> 47) on cells whose data
#     is no longer the most recently refreshed series ---
$ws.Range("N7").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("N22").PasteSpecial(-4122)
$ws.Range("N23").PasteSpecial(-4122)
$ws.Range("N37").PasteSpecial(-4122)
$ws.Range("N38").PasteSpecial(-4122)

# --- DTWEXBGS (row 39): now the most recently refreshed series -> add highlight (47 -> 48) ---
$ws.Range("N47").Copy() | Out-Null
$ws.Range("N39").PasteSpecial(-4122)
$ws.Range("N39").Value = 45989
$ws.Range("Q39").Value = 121.4288
$ws.Range("R39").ClearContents()
$ws.Range("S39").Value = 121.6225
$ws.Range("T39").Value = 122.0044
$ws.Range("U39").Value = 122.2833

# --- DFF (row 47) ---
$ws.Range("N47").Value = 45989
$ws.Range("Q47").Value = 3.89

# --- DGS2 (row 48) ---
$ws.Range("N48").Value = 45989
$ws.Range("Q48").Value = 3.47
$ws.Range("R48").ClearContents()
$ws.Range("S48").Value = 3.45
$ws.Range("T48").Value = 3.43
$ws.Range("U48").Value = 3.46

# --- DGS5 (row 49) ---
$ws.Range("N49").Value = 45989
$ws.Range("Q49").Value = 3.59
$ws.Range("R49").ClearContents()
$ws.Range("S49").Value = 3.56
$ws.Range("T49").Value = 3.55
$ws.Range("U49").Value = 3.61

# --- DGS10 (row 50) ---
$ws.Range("N50").Value = 45989
$ws.Range("Q50").Value = 4.02
$ws.Range("R50").ClearContents()
$ws.Range("S50").Value = 4
$ws.Range("T50").Value = 4.01
$ws.Range("U50").Value = 4.04

# --- DBAA (row 52) ---
$ws.Range("N52").Value = 45989
$ws.Range("Q52").Value = 5.8
$ws.Range("R52").ClearContents()
$ws.Range("S52").Value = 5.78
$ws.Range("T52").Value = 5.8
$ws.Range("U52").Value = 5.84
